$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Donor cell with default (unstyled) formatting, used to reset style after forcing text entry
$donorStyle = $ws.Cells.Item(45, 4).Style

$ws.Cells.Item(2, 4).Value = '93.318.18'
$ws.Cells.Item(2, 5).Value = '  +1.74%  '
$ws.Cells.Item(3, 4).Value = '3.096.54'
$ws.Cells.Item(3, 5).Value = '  -0.91%  '
$ws.Cells.Item(4, 5).Value = '  +0.00%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '236.62'
$ws.Cells.Item(5, 4).Style = $donorStyle
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '612.02'
$ws.Cells.Item(6, 4).Style = $donorStyle
$ws.Cells.Item(6, 5).Value = '  -1.03%  '
$ws.Cells.Item(7, 5).Value = '  +2.25%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.387'
$ws.Cells.Item(8, 4).Style = $donorStyle
$ws.Cells.Item(8, 5).Value = '  +0.80%  '
$ws.Cells.Item(9, 5).Value = '  -0.07%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.822'
$ws.Cells.Item(10, 4).Style = $donorStyle
$ws.Cells.Item(10, 5).Value = '  +11.47%  '
$ws.Cells.Item(11, 4).Value = '3.094.44'
$ws.Cells.Item(11, 5).Value = '  -0.87%  '
$ws.Cells.Item(12, 5).Value = '  -3.35%  '
$ws.Cells.Item(13, 5).Value = '  -3.18%  '
$ws.Cells.Item(14, 2).Value = 'WrappedBTC'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(14, 4).Value = '93.049.83'
$ws.Cells.Item(14, 5).Value = '  +1.62%  '
$ws.Cells.Item(15, 2).Value = 'Avalanche'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '34.77'
$ws.Cells.Item(15, 4).Style = $donorStyle
$ws.Cells.Item(15, 5).Value = '  -0.30%  '
$ws.Cells.Item(16, 5).Value = '  -3.36%  '
$ws.Cells.Item(17, 4).Value = '3.670.13'
$ws.Cells.Item(17, 5).Value = '  -0.94%  '
$ws.Cells.Item(18, 4).Value = '3.105.43'
$ws.Cells.Item(18, 5).Value = '  +0.52%  '
$ws.Cells.Item(19, 5).Value = '  -1.14%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '14.63'
$ws.Cells.Item(20, 4).Style = $donorStyle
$ws.Cells.Item(20, 5).Value = '  -2.19%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '5.96'
$ws.Cells.Item(21, 4).Style = $donorStyle
$ws.Cells.Item(21, 5).Value = '  +1.71%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '439.92'
$ws.Cells.Item(22, 4).Style = $donorStyle
$ws.Cells.Item(22, 5).Value = '  -2.02%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '0.0000198'
$ws.Cells.Item(23, 4).Style = $donorStyle
$ws.Cells.Item(23, 5).Value = '  -1.33%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '8.99'
$ws.Cells.Item(24, 4).Style = $donorStyle
$ws.Cells.Item(24, 5).Value = '  -5.33%  '
$ws.Cells.Item(25, 5).Value = '  +4.54%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '5.66'
$ws.Cells.Item(26, 4).Style = $donorStyle
$ws.Cells.Item(26, 5).Value = '  -3.96%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '12.79'
$ws.Cells.Item(27, 4).Style = $donorStyle
$ws.Cells.Item(27, 5).Value = '  +8.53%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '85.76'
$ws.Cells.Item(28, 4).Style = $donorStyle
$ws.Cells.Item(28, 5).Value = '  -2.59%  '
$ws.Cells.Item(29, 5).Value = '  -0.21%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '0.251'
$ws.Cells.Item(30, 4).Style = $donorStyle
$ws.Cells.Item(30, 5).Value = '  +5.81%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '0.181'
$ws.Cells.Item(31, 4).Style = $donorStyle
$ws.Cells.Item(31, 5).Value = '  +7.72%  '
$ws.Cells.Item(32, 5).Value = '  -16.03%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '9.18'
$ws.Cells.Item(33, 4).Style = $donorStyle
$ws.Cells.Item(33, 5).Value = '  -2.31%  '
$ws.Cells.Item(34, 5).Value = '  +0.68%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '7.93'
$ws.Cells.Item(35, 4).Style = $donorStyle
$ws.Cells.Item(35, 5).Value = '  +0.36%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.158'
$ws.Cells.Item(36, 4).Style = $donorStyle
$ws.Cells.Item(36, 5).Value = '  -10.07%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '25.84'
$ws.Cells.Item(37, 4).Style = $donorStyle
$ws.Cells.Item(37, 5).Value = '  -1.73%  '
$ws.Cells.Item(38, 5).Value = '  -5.86%  '
$ws.Cells.Item(39, 5).Value = '  -2.58%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.451'
$ws.Cells.Item(40, 4).Style = $donorStyle
$ws.Cells.Item(40, 5).Value = '  +1.29%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '23.98'
$ws.Cells.Item(41, 4).Style = $donorStyle
$ws.Cells.Item(41, 5).Value = '  +8.09%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '473.93'
$ws.Cells.Item(42, 4).Style = $donorStyle
$ws.Cells.Item(43, 5).Value = '  -2.16%  '
$ws.Cells.Item(44, 5).Value = '  -4.13%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '159.16'
$ws.Cells.Item(46, 4).Style = $donorStyle
$ws.Cells.Item(46, 5).Value = '  +0.29%  '
$ws.Cells.Item(47, 5).Value = '  -1.54%  '
$ws.Cells.Item(48, 5).Value = '  -3.45%  '
$ws.Cells.Item(49, 2).Value = 'OKB'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '43.77'
$ws.Cells.Item(49, 4).Style = $donorStyle
$ws.Cells.Item(49, 5).Value = '  -0.62%  '
$ws.Cells.Item(50, 2).Value = 'ImmutableX'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '1.31'
$ws.Cells.Item(50, 4).Style = $donorStyle
$ws.Cells.Item(50, 5).Value = '  -3.81%  '
$ws.Cells.Item(51, 2).Value = 'VeChain'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.0323'
$ws.Cells.Item(51, 4).Style = $donorStyle
$ws.Cells.Item(51, 5).Value = '  -1.06%  '
